$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the anchor paragraphs by their (stable, unique) text so this script
# does not depend on hard-coded paragraph indices.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$startIndex = -1
$endIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -match "^Size = ") {
        $startIndex = $i
    }
    if ($t -match "^====") {
        $endIndex = $i
    }
}

if ($startIndex -eq -1 -or $endIndex -eq -1) {
    throw "Could not locate anchor paragraphs (start=$startIndex end=$endIndex)"
}

$startPara = $paras.Item($startIndex)
$endPara = $paras.Item($endIndex)

$target = $d.Range($startPara.Range.Start, $endPara.Range.End)

# ---------------------------------------------------------------------------
# Replacement content (OOXML) covering everything from the old "Size = ..."
# paragraph through the old "====...====" paragraph. This both rewrites
# existing paragraphs (formatting / text tweaks) and inserts the new
# "Expected changes" / "Initialise Board" sections that now precede the
# re-ordered determinator-logic write-up.
# ---------------------------------------------------------------------------
$newXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
'<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>Expected changes</w:t></w:r></w:p>' + `
'<w:p><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Number of consecutive </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>coords</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> needs to be determined. Will start with 3 as per the basic game.</w:t></w:r></w:p>' + `
'<w:p/>' + `
'<w:p><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Initialise Board</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Write to console \u201cWhat size would you like the board?\u201d</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Read input from console</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Send input to Print method</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Print method writes board to console</w:t></w:r></w:p>' + `
'<w:p/>' + `
'<w:p><w:r><w:t xml:space="preserve">Size = </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>board.length</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>NB determined by players \u2013 from 3x3 to 10x10</w:t></w:r><w:r><w:t xml:space="preserve"> (defined as const)</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Reverse = size \u2013 1</w:t></w:r><w:r><w:t xml:space="preserve"> \u2013 defined on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>DiagonalDeterminator</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' + `
'<w:p><w:r><w:t>** Replace for loop through rows with foreach loop</w:t></w:r></w:p>' + `
'<w:p/>' + `
'<w:p><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Column determinator logic:</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Count == 0</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>*</w:t></w:r><w:r><w:t>*</w:t></w:r><w:r><w:t>For loop through rows</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>For loop through columns</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t xml:space="preserve">If board[col][row] equals </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>player.PlayerId</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p>' + `
'<w:p><w:r><w:t>Count++</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>If count equals 3</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Return true</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Break</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Else</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Count = 0</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Return false</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>End loop through columns</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>End loop through rows</w:t></w:r></w:p>' + `
'<w:p/>' + `
'<w:p><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Row determinator logic:</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Count == 0</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>*</w:t></w:r><w:r><w:t>*</w:t></w:r><w:r><w:t>For loop through rows</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>For loop through columns</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t xml:space="preserve">If board[row][column] equals </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>player.PlayerId</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p>' + `
'<w:p><w:r><w:t>Count++</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>If count equals 3</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Return true</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Break</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Else</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Count = 0</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Return false</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>End loop through columns</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>End loop through rows</w:t></w:r></w:p>' + `
'<w:p/>' + `
'<w:p><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Diagonal determinator logic:</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Count == 0</w:t></w:r></w:p>' + `
'<w:p><w:r><w:lastRenderedPageBreak/><w:t>*</w:t></w:r><w:r><w:t>*</w:t></w:r><w:r><w:t>For loop through rows</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>For loop through columns</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t xml:space="preserve">If board[row][row] equals </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>player.PlayerId</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p>' + `
'<w:p><w:r><w:t>Or</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t xml:space="preserve">If board[row][reverse] equals </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>player.PlayerId</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p>' + `
'<w:p><w:r><w:t>Count++</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>If count equals 3</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Return true</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Break</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Else</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Count = 0</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>Return false</w:t></w:r></w:p>' + `
'<w:p><w:r><w:t>End loop through columns</w:t></w:r></w:p>' + `
'<w:p><w:pPr><w:pBdr><w:bottom w:val="double" w:sz="6" w:space="1" w:color="auto"/></w:pBdr></w:pPr><w:r><w:t>End loop through rows</w:t></w:r></w:p>' + `
'<w:p/>' + `
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($newXml)

Write-Output "Done. Paragraph count now: $($d.Paragraphs.Count)"
